# "running D suite only"
# Sheet 1 ("Test Cases") contains B-suite test cases in rows 2-28.
# Rows 16-24 correspond to TestCase_B15..B23: turn off their Runmode (column C) from "Y" to "N".
# Row 26 corresponds to TestCase_B25: mark its Results (column D) as "FAIL" instead of "SKIP".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

foreach ($r in 16..24) {
    $ws.Cells.Item($r, 3).Value = "N"
}

$ws.Cells.Item(26, 4).Value = "FAIL"

# Update the active selection to reflect the last edited cell.
$ws.Activate()
$ws.Range("B28").Select()
